$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "PAR (DNN)"
$ws.Range("A5").Value = "Temperature (DNN)"
$ws.Range("A8").Value = "Humidity (DNN)"
$ws.Range("A11").Value = "CO2 (DNN)"
$ws.Range("A14").Value = "Leaf Temperature (DNN)"
